$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 data (foldchange values), with A2 styled like the header row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.4951257660418741
$ws.Range("C2").Value = 0.2990535882639475
$ws.Range("D2").Value = 0.6825846768843439
$ws.Range("E2").Value = 0.3458615801392663
$ws.Range("F2").Value = 1.552836929124044
$ws.Range("G2").Value = 1.214498745700118

# Copy the header-row formatting (bold/border/centered) onto A2
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
